$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 10
$ws.Range("K4").Value = 15
$ws.Range("K5").Value = 20
$ws.Range("K6").Value = 25

$ws.Range("K8").Value = 10
$ws.Range("K9").Value = 15
$ws.Range("K10").Value = 20
$ws.Range("K11").Value = 25

$ws.Range("K13").Value = 10
$ws.Range("K14").Value = 15
$ws.Range("K15").Value = 20
$ws.Range("K16").Value = 25
